$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text (e.g. "148.98", "0.0490") that must stay as
# literal text, not be auto-coerced to a number (which would drop trailing
# zeros / change formatting). Format those cells as Text before assigning.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '71.621.28'
$ws.Range('E2').Value = '  +2.63%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '4.028.00'
$ws.Range('E3').Value = '  +2.22%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('E5').Value = '  -1.07%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '148.98'
$ws.Range('E6').Value = '  +1.62%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.625'
$ws.Range('E7').Value = '  +0.45%  '
$ws.Range('E8').Value = '  +0.24%  '
$ws.Range('E9').Value = '  +0.96%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.176'
$ws.Range('E10').Value = '  +1.39%  '
$ws.Range('E11').Value = '  +0.00%  '
$ws.Range('E12').Value = '  +6.75%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '10.80'
$ws.Range('E13').Value = '  +3.02%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.668.00'
$ws.Range('E14').Value = '  +2.20%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.044.73'
$ws.Range('E15').Value = '  +2.83%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '21.40'
$ws.Range('E16').Value = '  +7.40%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '14.31'
$ws.Range('E17').Value = '  +1.49%  '
$ws.Range('E18').Value = '  +0.88%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.133'
$ws.Range('E19').Value = '  -1.79%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '71.606.13'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '442.07'
$ws.Range('E21').Value = '  +1.55%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '3.61'
$ws.Range('E22').Value = '  +5.81%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '95.10'
$ws.Range('E23').Value = '  +7.29%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '14.36'
$ws.Range('E24').Value = '  -1.67%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '12.24'
$ws.Range('E25').Value = '  +2.51%  '
$ws.Range('E26').Value = '  -1.67%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '11.18'
$ws.Range('E27').Value = '  +1.62%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '37.10'
$ws.Range('E28').Value = '  +0.98%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '13.62'
$ws.Range('E29').Value = '  +1.67%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '699.50'
$ws.Range('E30').Value = '  -1.12%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.132'
$ws.Range('E31').Value = '  +2.73%  '
$ws.Range('E32').Value = '  +1.65%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '6.95'
$ws.Range('E33').Value = '  +14.23%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '67.92'
$ws.Range('E34').Value = '  -1.08%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0₃0906'
$ws.Range('E35').Value = '  +3.87%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.444'
$ws.Range('E36').Value = '  -0.33%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '41.03'
$ws.Range('E37').Value = '  +0.78%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.159'
$ws.Range('E38').Value = '  +6.06%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.58'
$ws.Range('E39').Value = '  +18.77%  '
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('E41').Value = '  -0.27%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0490'
$ws.Range('E42').Value = '  +0.94%  '
$ws.Range('E43').Value = '  +0.37%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '3.12'
$ws.Range('E44').Value = '  +0.79%  '
$ws.Range('E45').Value = '  +3.31%  '
$ws.Range('E46').Value = '  +2.26%  '
$ws.Range('E47').Value = '  -3.54%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '9.23'
$ws.Range('E48').Value = '  +5.94%  '
$ws.Range('E49').Value = '  +0.74%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.000276'
$ws.Range('E50').Value = '  +16.91%  '

# Row 51: coin swapped from BabyDogeCoin to ARBITRUM
$ws.Range('B51').Value = 'ARBITRUM'
$ws.Range('C51').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.08'
$ws.Range('E51').Value = '  -0.90%  '
